# Fix typo "Tripple Residential Pod" -> "Triple Residential Pod".
#
# The canonical name lives in the "Property Types" table (cell B12), but
# the "Properties" sheet stores the Property Type as literal text in its
# own column C (it is the VLOOKUP *key*, not a formula result), so every
# occurrence there must be corrected too or the lookups on the
# Properties/Leases sheets will fail to resolve against the renamed key.

$wb = $excel.ActiveWorkbook

$wsTypes = $wb.Worksheets.Item("Property Types")
$wsTypes.Range("B12").Value = "Triple Residential Pod"

$wsProperties = $wb.Worksheets.Item("Properties")
for ($r = 11; $r -le 60; $r++) {
    $cell = $wsProperties.Cells.Item($r, 3)
    if ($cell.Text -eq "Tripple Residential Pod") {
        $cell.Value = "Triple Residential Pod"
    }
}

# Recalculate so the dependent VLOOKUP formulas on the Properties and
# Leases sheets pick up the corrected text.
$excel.Calculate()

# Update the active selections on each sheet to match the author's final
# cursor positions.
$wsLeases = $wb.Worksheets.Item("Leases")
$wsLeases.Activate() | Out-Null
$wsLeases.Range("A6").Select() | Out-Null

$wsProperties = $wb.Worksheets.Item("Properties")
$wsProperties.Activate() | Out-Null
$wsProperties.Range("A2").Select() | Out-Null

$wsTypes.Activate() | Out-Null
$wsTypes.Range("B11").Select() | Out-Null

$wsLeases.Activate() | Out-Null
